$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6441636085510254
$ws.Range("B1").Value = 1.990078449249268
$ws.Range("C1").Value = 6.531002521514893
$ws.Range("D1").Value = 1.583864092826843
$ws.Range("E1").Value = 0.9095044136047363
